# Insert a new weekly price record for "Feria Lagunitas de Puerto Montt - Albahaca"
# as row 66, pushing the existing rows 66:112 down to 67:113.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 66; this shifts rows 66-112 to 67-113
# and carries the existing row 66 formatting (incl. the date number format on column D)
# down with it.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new record.
$ws.Cells.Item(66, 1).Value = 4
$ws.Cells.Item(66, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(66, 3).Value = "Los Lagos"
$ws.Cells.Item(66, 4).Value = 44651
$ws.Cells.Item(66, 5).Value = 10
$ws.Cells.Item(66, 6).Value = 100112052
$ws.Cells.Item(66, 7).Value = "Albahaca"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 40
$ws.Cells.Item(66, 11).Value = 6000
$ws.Cells.Item(66, 12).Value = 6000
$ws.Cells.Item(66, 13).Value = 6000
$ws.Cells.Item(66, 14).Value = "`$/docena de matas"
$ws.Cells.Item(66, 15).Value = "Región Metropolitana"
$ws.Cells.Item(66, 16).Value = 1000
$ws.Cells.Item(66, 17).Value = 6
$ws.Cells.Item(66, 18).Value = "Hortaliza"
